$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$data = @(
  @("2022-06", 0.000334165, 0.000329191, 0.0002679798819084184, 0.9293689819276465, 0.9293896957646928),
  @("2022-07", 0.000316309, 0.000313162, 0.0002546446245606605, 0.9629654914481928, 0.9629749416417263),
  @("2022-08", 0.000307203, 0.000305671, 0.0002482733755398285, 1.009369342244816, 1.009370346514574),
  @("2022-09", 0.000298527, 0.000298527, 0.0002421977174380582, 1.046125365123627, 1.046120104947815),
  @("2022-10", 0.000290721, 0.000292178, 0.0002367797544677461, 1.075813196343176, 1.075803585747267),
  @("2022-11", 0.000258778, 0.000261379, 0.0002115812832840122, 1.002842220523855, 1.002830937019458),
  @("2022-12", 0.000277715, 0.000281912, 0.0002279454696919372, 1.119231652089665, 1.119217776450041),
  @("2022-13", 0.000269452, 0.000274895, 0.0002220213668001376, 1.122739911357767, 1.12272566034234),
  @("2022-14", 0.000249327, 0.000255638, 0.0002062359657351289, 1.068851528790272, 1.068838289280802),
  @("2022-15", 0.000226992, 0.000233904, 0.000188489593315789, 0.9970729905278174, 0.9970613562440914),
  @("2022-16", 0.000216062, 0.000223758, 0.0001801096368916599, 0.9691031890630956, 0.9690928364147705),
  @("2022-17", 0.000210837, 0.000219441, 0.0001764359660070377, 0.9628588993431316, 0.9628497157420473),
  @("2022-18", 0.000227141, 0.000237594, 0.0001908172171375343, 1.053628040524997, 1.053619280548596),
  @("2022-19", 0.000217962, 0.000229136, 0.0001838165207031035, 1.024873465120994, 1.024866219359078),
  @("2022-20", 0.000180642, 0.000190855, 0.0001529340617102859, 0.8595476836421185, 0.8595426580354665),
  @("2022-21", 0.000171881, 0.000182509, 0.0001460814579019376, 0.8264646668887552, 0.8264608044807564),
  @("2022-22", 0.000178947, 0.000190964, 0.0001526769155358249, 0.8684591258480909, 0.8684560267416269),
  @("2022-23", 0.000194813, 0.000208937, 0.0001668586104604379, 0.9533196758569621, 0.9533172510904792),
  @("2022-24", 0.000167136, 0.000180152, 0.0001437084897294734, 0.823993119213801, 0.8239917973110351),
  @("2022-25", 0.000182124, 0.000197291, 0.0001572031825854582, 0.9039612578183586, 0.9039605774099464),
  @("2022-26", 0.00020064, 0.000218439, 0.0001738574901974597, 1.00201582274271, 1.002015834953403),
  @("2022-27", 0.00017163, 0.000187791, 0.0001492969109543221, 0.8620080715739403, 0.8620086692468327),
  @("2022-28", 0.000179143, 0.000196994, 0.0001564368977855694, 0.9044827273569107, 0.9044838987310392),
  @("2022-29", 0.000197228, 0.000217969, 0.0001728978750948279, 1.000694532318924, 1.000696356045175),
  @("2022-30", 0.000206956, 0.000229866, 0.0001821297283045357, 1.054916288500898, 1.054918694798929),
  @("2022-31", 0.000192022, 0.000214348, 0.0001696428419695281, 0.9830889852890349, 0.9830916165091936),
  @("2022-32", 0.000200871, 0.000225349, 0.0001781490558780481, 1.032696822522912, 1.032699935550818),
  @("2022-33", 0.000189894, 0.000214102, 0.0001690671606147519, 0.9801821195587758, 0.9801853546095168),
  @("2022-34", 0.00018993, 0.000215216, 0.0001697552901657529, 0.9841633316349834, 0.9841668152567812),
  @("2022-35", 0.000163076, 0.000185713, 0.0001463192721718276, 0.8481836819244233, 0.8481868515976801),
  @("2022-36", 0.00021513, 0.000246221, 0.0001937734144014444, 1.123010107921022, 1.1230144842857),
  @("2022-37", 0.000188274, 0.000216564, 0.0001702414450354405, 0.9863237203990848, 0.9863276891577741),
  @("2022-38", 0.000206394, 0.000238597, 0.0001873500218791518, 1.085031007949167, 1.085035479569404),
  @("2022-39", 0.000218791, 0.000254195, 0.0001993736989164537, 1.154157144356778, 1.15416198332018),
  @("2022-40", 0.000229871, 0.000268407, 0.0002102830807037777, 1.216714496326974, 1.216719656687607),
  @("2022-41", 0.000227717, 0.000267224, 0.0002091208484473001, 1.209347257460768, 1.20935242104328),
  @("2022-42", 0.000195981, 0.000231135, 0.0001806747992006032, 1.044252652109058, 1.044257122155176),
  @("2022-43", 0.000205734, 0.000243854, 0.000190401959755907, 1.099818814467124, 1.099823517399889),
  @("2022-44", 0.000185902, 0.000221451, 0.0001727154303512224, 0.997038767948026, 0.9970430134564792),
  @("2022-45", 0.000194771, 0.000233179, 0.0001816574040990949, 1.047988089779116, 1.047992521091982),
  @("2022-46", 0.000185972, 0.000223762, 0.0001741237870665303, 1.00386723237584, 1.003871437036129),
  @("2022-47", 0.000186891, 0.000225994, 0.0001756631487384419, 1.012063036385649, 1.012067225961724),
  @("2022-48", 0.000200627, 0.00024382, 0.0001893055677033155, 1.089917972941813, 1.08992242318698),
  @("2022-49", 0.00020951, 0.000255891, 0.0001984542928591538, 1.141800148466962, 1.141804738356822),
  @("2022-50", 0.000215744, 0.000264826, 0.0002051522102154817, 1.179509159903074, 1.17951381996198),
  @("2022-51", 0.000242329, 0.00029895, 0.0002313260790399442, 1.329052621226495, 1.329057773704098),
  @("2022-52", 0.000270703, 0.000335628, 0.000259414357056701, 1.489366177847251, 1.489371835097677),
  @("2023-01", 0.000222539, 0.000277295, 0.0002140862644250865, 1.228242259236588, 1.228246823762069),
  @("2023-02", 0.000206655, 0.000258793, 0.0001995769229005035, 1.14417250483532, 1.144176659274492),
  @("2023-03", 0.000198286, 0.000249558, 0.0001922375352409787, 1.10129538779, 1.101299289479022),
  @("2023-04", 0.000183271, 0.000231816, 0.0001783699210149897, 1.021105073063005, 1.02110859820464),
  @("2023-05", 0.000189062, 0.00024034, 0.0001847199707057014, 1.056683168260688, 1.05668671837641),
  @("2023-06", 0.000182011, 0.000232536, 0.000178520860606054, 1.02047213068303, 1.02047546282985),
  @("2023-07", 0.000186031, 0.000238863, 0.0001831717056315018, 1.046287559039711, 1.046290875135635),
  @("2023-08", 0.000175874, 0.000226954, 0.0001738426936797218, 0.9922680033065064, 0.9922710517260998),
  @("2023-09", 0.000182553, 0.000236753, 0.0001811446486847976, 1.033183363511029, 1.033186436020881),
  @("2023-10", 0.000178597, 0.000232783, 0.0001779067502698257, 1.013965635562871, 1.013968550223506),
  @("2023-11", 0.000173752, 0.000227604, 0.000173752, 0.9895532512298802, 0.9895559966218566),
  @("2023-12", 0.000176886, 0.00023287, 0.0001775722893824233, 1.010561602797019, 1.010564304557344),
  @("2023-13", 0.000176474, 0.000233492, 0.0001778460382677795, 1.011369326977883, 1.011371928283884),
  @("2023-14", 0.000174287, 0.000231755, 0.0001763234979232004, 1.001967291250911, 1.001969766192165),
  @("2023-15", 0.000162784, 0.000217543, 0.0001653250440595994, 0.938770830875166, 0.9387730535599438),
  @("2023-16", 0.000147725, 0.000198407, 0.0001506130705606076, 0.8545963613714278, 0.8545982968997444),
  @("2023-17", 0.000153959, 0.000207817, 0.0001575779616399733, 0.8934518766731073, 0.8934538080564249),
  @("2023-18", 0.000162415, 0.000220331, 0.000166877684845185, 0.9454774107480255, 0.9454793568192563),
  @("2023-19", 0.000132701, 0.000180923, 0.0001368762372708105, 0.7749221313829677, 0.7749236460781647),
  @("2023-20", 0.000148257, 0.000203144, 0.0001535149950692156, 0.8684760572208017, 0.868477664583452),
  @("2023-21", 0.000144283, 0.00019869, 0.0001499797038508069, 0.8478451045975999, 0.8478465855750932),
  @("2023-22", 0.000138086, 0.00019111, 0.0001440949332261887, 0.8139723662883067, 0.8139737032997554),
  @("2023-23", 0.000148321, 0.000206303, 0.0001553758216893204, 0.8770437772268783, 0.8770451263574913),
  @("2023-24", 0.000143012, 0.000199916, 0.0001503955570139484, 0.8483004257269555, 0.8483016420685772),
  @("2023-25", 0.000155028, 0.000217799, 0.0001636644680851098, 0.9224565019301713, 0.9224577281992424),
  @("2023-26", 0.000141277, 0.000199476, 0.0001497260805525879, 0.8432682281813552, 0.8432692609890114),
  @("2023-27", 0.000148852, 0.000211225, 0.0001583661646269472, 0.8912662224577647, 0.8912672207800224),
  @("2023-28", 0.000155986, 0.000222457, 0.0001666000313955729, 0.9369078479067565, 0.9369087992513458),
  @("2023-29", 0.000138229, 0.000198121, 0.0001482075596781807, 0.8328540235232674, 0.8328547819638896),
  @("2023-30", 0.00013647, 0.00019658, 0.0001468892842313463, 0.8248317456385673, 0.8248324103331941),
  @("2023-31", 0.00012804, 0.000185362, 0.0001383503688395348, 0.7763048236068042, 0.776305367813487),
  @("2023-32", 0.000149402, 0.000217372, 0.0001620588675983001, 0.9086601138443445, 0.9086606555503038),
  @("2023-33", 0.000175667, 0.000256867, 0.0001912882552649019, 1.071750426598994, 1.071750953118626),
  @("2023-34", 0.000156569, 0.000230088, 0.0001711534398228618, 0.9582254004302045, 0.9582257706509055),
  @("2023-35", 0.000126338, 0.000186593, 0.0001386422457336413, 0.7756294918795633, 0.7756297101668023),
  @("2023-36", 0.000139704, 0.000207366, 0.0001539047985518737, 0.8603745758282142, 0.8603747276742519),
  @("2023-37", 0.000138833, 0.000207107, 0.0001535386652560339, 0.8576890055921148, 0.8576890669433838),
  @("2023-38", 0.000147754, 0.00022152, 0.0001640385917654711, 0.9156611452055052, 0.9156611145878115),
  @("2023-39", 0.000153118, 0.000230713, 0.0001706533295127182, 0.9518755867191163, 0.9518754549641955),
  @("2023-40", 0.000166499, 0.000252132, 0.0001862867132984335, 1.038302666858229, 1.03830241413266),
  @("2023-41", 0.000148714, 0.000226328, 0.0001670335990343469, 0.9302990828159661, 0.9302987587034688),
  @("2023-42", 0.00016477, 0.000252021, 0.0001857855182384851, 1.033968422319676, 1.033967953523227),
  @("2023-43", 0.00015277, 0.000234837, 0.0001729233052183204, 0.9616689794808944, 0.9616684424852928),
  @("2023-44", 0.000160813, 0.00024844, 0.0001827335687469489, 1.01546995015161, 1.015469276480541),
  @("2023-45", 0.000170642, 0.000264947, 0.0001946546789978841, 1.080911774947295, 1.080910944352718),
  @("2023-46", 0.000177802, 0.000277448, 0.0002036091493559278, 1.129794221674025, 1.129793234871616),
  @("2023-47", 0.00016669, 0.000261411, 0.000191624893757869, 1.062504175601952, 1.062503135991222),
  @("2023-48", 0.00017073, 0.000269089, 0.0001970307241930985, 1.091664858320368, 1.091663675530646),
  @("2023-49", 0.000192164, 0.000304389, 0.0002226270216697089, 1.232565225117328, 1.232563760219551),
  @("2023-50", 0.00016767, 0.000266922, 0.0001950037410646629, 1.078826572273562, 1.07882517679137),
  @("2023-51", 0.000200263, 0.000320406, 0.0002338137403358645, 1.292573887705998, 1.292572079985432),
  @("2023-52", 0.000191379, 0.000307727, 0.0002243082904367699, 1.239102778923391, 1.239100915846623),
  @("2024-01", 0.000196771, 0.000317983, 0.0002315228553763729, 1.278004938720912, 1.278002882927233),
  @("2024-02", 0.000159763, 0.000259471, 0.0001887081760809532, 1.040892874635776, 1.040891090937091),
  @("2024-03", 0.00017452, 0.000284859, 0.0002069385716501475, 1.140600143598363, 1.140598069243609),
  @("2024-04", 0.000176783, 0.000289999, 0.0002104352427174218, 1.159009794538149, 1.159007564973109),
  @("2024-05", 0.000163863, 0.000270153, 0.000195812589195232, 1.077670257495002, 1.077668071214219),
  @("2024-06", 0.00017729, 0.000293753, 0.0002126795247665827, 1.169627643533715, 1.169625147852261),
  @("2024-07", 0.000159005, 0.000264778, 0.0001914846455184996, 1.052282998956895, 1.05228064313697),
  @("2024-08", 0.000149648, 0.000250445, 0.0001809155195748897, 0.9934615843389831, 0.9934592558629719),
  @("2024-09", 0.000142521, 0.000239713, 0.0001729678951421364, 0.949111882669477, 0.9491095584544484),
  @("2024-10", 0.000148798, 0.000251525, 0.0001812865007839271, 0.9940175040712382, 0.9940149654873555),
  @("2024-11", 0.000149714, 0.000254342, 0.0001831101924076832, 1.00326980238384, 1.003267134796352),
  @("2024-12", 0.000130514, 0.000222836, 0.0001602466421983244, 0.877345828044497, 0.8773434031267321),
  @("2024-13", 0.000144838, 0.000248531, 0.000178523786612721, 0.9766852358452927, 0.9766824337783275),
  @("2024-14", 0.000131445, 0.000226679, 0.0001626444972087718, 0.8891490373116074, 0.8891463929935355),
  @("2024-15", 0.000142195, 0.000246447, 0.0001766287363520828, 0.9648797807716721, 0.9648768098885238),
  @("2024-16", 0.000146241, 0.00025473, 0.0001823593003212859, 0.9954430451230957, 0.9954398755821035)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Remove the now-unused last row (117), since new dataset has only 115 data rows (2..116)
$ws.Rows.Item(117).Delete()

# Update parameters sheet
$ws2 = $wb.Worksheets.Item("parameters")
$ws2.Cells.Item(2, 3).Value = -8.63080551131233
$ws2.Cells.Item(3, 3).Value = -0.003872332855820278
$ws2.Cells.Item(4, 3).Value = 57
$ws2.Cells.Item(5, 3).Value = -7.9306262378093
$ws2.Cells.Item(6, 3).Value = -0.003127804535891093
$ws2.Cells.Item(7, 3).Value = -0.09837518444799552
$ws2.Cells.Item(8, 3).Value = 5.653607367450982
$ws2.Cells.Item(9, 3).Value = -8.05022281338894
$ws2.Cells.Item(10, 3).Value = -0.05491529906278651
$ws2.Cells.Item(11, 3).Value = 0.1082490153532271
$ws2.Cells.Item(12, 3).Value = 38
$ws2.Cells.Item(17, 3).Value = -7.930648525626751
$ws2.Cells.Item(18, 3).Value = -0.003127699503525955
$ws2.Cells.Item(19, 3).Value = -0.09836054126441914
$ws2.Cells.Item(20, 3).Value = 5.65432870282719
$ws2.Cells.Item(21, 3).Value = -8.05022281338894
$ws2.Cells.Item(22, 3).Value = -0.05491529906278651
$ws2.Cells.Item(23, 3).Value = 0.1082490153532271
$ws2.Cells.Item(24, 3).Value = 38
